# This script reproduces a manual Excel edit where the user selected the
# range G5:G11 (col G, "Fonte FINEP" source-column data) and deleted the
# cell G5 ("Data da liberação"), shifting the remaining cells in column G
# up by one row (G6->G5, G7->G6, ... G12->G11), leaving G12 empty.
#
# Net effect on sheet "Planilha1":
#   G5:  Data da liberação   -> Data da assinatura
#   G6:  Data da assinatura  -> Prazo de utilização
#   G7:  Prazo de utilização -> Contrato
#   G8:  Contrato            -> Instrumento
#   G9:  Instrumento         -> Proponente
#   G10: Proponente          -> UF do Executor
#   G11: UF do Executor      -> Status do projeto
#   G12: Status do projeto   -> (empty)
#
# Columns A-F and rows 1-4 / 13-15 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Reproduce the selection left behind by the edit (matches the saved
# sheetView/selection: activeCell="G5" sqref="G5:G11").
$null = $ws.Range("G5:G11").Select()

# Shift column G up by one row across rows 5..11 (copy value from the row
# below into the current row), then clear the now-vacated last row (G12).
for ($r = 5; $r -le 11; $r++) {
    $below = $ws.Cells.Item($r + 1, 7)
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = $below.Value2
}
$ws.Cells.Item(12, 7).ClearContents()
